# Add a new row (row 4) to the "CLIENTES_PERMISOS" sheet with a new client:
#   A4 = 946 (CUENTA, plain number, default style)
#   B4 = "FABIAN FERRETERIA" (NOMBRE)
#   C4 = "SANITARIOS,REPUESTOS,BULONERIA,PINTURAS,FERRETERIA,ELECTRICIDAD" (CATEGORIAS)
#   D4 = "D" (LISTA_PRECIOS)
# and move the active selection to D5, just below the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row's values.
$ws.Cells.Item(4, 1).Value = 946
$ws.Cells.Item(4, 2).Value = "FABIAN FERRETERIA"
$ws.Cells.Item(4, 3).Value = "SANITARIOS,REPUESTOS,BULONERIA,PINTURAS,FERRETERIA,ELECTRICIDAD"
$ws.Cells.Item(4, 4).Value = "D"

# Copy the formatting (font/style) from the equivalent cells of row 2 (B2:D2)
# onto the new row's B4:D4 cells, matching how the other data rows look.
# Column A is left with the default style, matching the source edit.
$ws.Range("B2:D2").Copy()
$ws.Range("B4:D4").PasteSpecial(-4122)

# Match the row height used by the rest of the data rows.
$ws.Rows(4).RowHeight = 15.75

# Move the selection to D5, as in the edited workbook.
$ws.Range("D5").Select()
